# Update countries & provincias Spain
#
# - Reorder the "Ceuta" / "Lanzarote" entries (Ceuta now listed before
#   Lanzarote) and refresh Ceuta's case counters.
# - Bump the "last updated" timestamp in the header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58 becomes Ceuta with its updated figures.
$ws.Range("A58").Value = "Ceuta"
$ws.Range("B58").Value = 20
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 19
$ws.Range("E58").Value = 1

# Row 59 becomes Lanzarote, keeping its previous figures.
$ws.Range("A59").Value = "Lanzarote"
$ws.Range("B59").Value = 17
$ws.Range("C59").Value = 15
$ws.Range("D59").Value = 17
$ws.Range("E59").Value = 36

# Refresh the "datos actualizados" timestamp (15:29 -> 15:59).
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 15:59"
